$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only A2 (signal rate) changes; B2..H2 remain as-is
$ws.Range("A2").Value = 0.00036785000702366233

# Row 3
$ws.Range("A3").Value = 0.00010025110532296821
$ws.Range("B3").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_state7"
$ws.Range("C3").Value = 36.349205017089844
$ws.Range("D3").Value = 10.912285804748535
$ws.Range("E3").Value = 88.0
$ws.Range("F3").Value = 28.0
$ws.Range("G3").Value = "myclk"
$ws.Range("H3").Value = "DSP FF LUT "

# Row 4
$ws.Range("A4").Value = 0.00009444444731343538
$ws.Range("B4").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_state13"
$ws.Range("C4").Value = 39.68254089355469
$ws.Range("D4").Value = 11.904762268066406
$ws.Range("E4").Value = 105.0
$ws.Range("F4").Value = 21.0
$ws.Range("G4").Value = "myclk"
$ws.Range("H4").Value = "FF LUT "

# Row 5
$ws.Range("A5").Value = 0.00003508024019538425
$ws.Range("B5").Value = "firConvolutionOperationChaining_IP/U0/shiftRegister_U/firConvolutionOpebkb_ram_U/q0[31]_i_1_n_0"
$ws.Range("C5").Value = 39.841270446777344
$ws.Range("D5").Value = 21.577285766601562
$ws.Range("E5").Value = 32.0
$ws.Range("F5").Value = 14.0
$ws.Range("G5").Value = "myclk"
$ws.Range("H5").Value = "FF "

# Row 6
$ws.Range("A6").Value = 0.000027975236662314273
$ws.Range("B6").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_state12"
$ws.Range("C6").Value = 36.19047546386719
$ws.Range("D6").Value = 10.85714340209961
$ws.Range("E6").Value = 34.0
$ws.Range("F6").Value = 9.0
$ws.Range("G6").Value = "myclk"
$ws.Range("H6").Value = "FF LUT "

# Row 7
$ws.Range("A7").Value = 0.000026315159630030394
$ws.Range("B7").Value = "firConvolutionOperationChaining_IP/U0/shiftRegister_U/firConvolutionOpebkb_ram_U/p_0_in_0"
$ws.Range("C7").Value = 39.841270446777344
$ws.Range("D7").Value = 11.934745788574219
$ws.Range("E7").Value = 32.0
$ws.Range("F7").Value = 8.0
$ws.Range("G7").Value = "myclk"
$ws.Range("H7").Value = "RAM "

# Row 8
$ws.Range("A8").Value = 0.00002261904774059076
$ws.Range("B8").Value = "firConvolutionOperationChaining_IP/U0/p_pn_reg_118[31]_i_1_n_0"
$ws.Range("C8").Value = 39.68254089355469
$ws.Range("D8").Value = 11.904762268066406
$ws.Range("E8").Value = 32.0
$ws.Range("F8").Value = 10.0
$ws.Range("G8").Value = "myclk"
$ws.Range("H8").Value = "FF "

# Row 9
$ws.Range("A9").Value = 0.000021919999198871665
$ws.Range("B9").Value = "firConvolutionOperationChaining_IP/U0/ce0"
$ws.Range("C9").Value = 43.4920654296875
$ws.Range("D9").Value = 13.047618865966797
$ws.Range("E9").Value = 17.0
$ws.Range("F9").Value = 8.0
$ws.Range("G9").Value = "myclk"
$ws.Range("H9").Value = "FF LUT "

# Row 10
$ws.Range("A10").Value = 0.000009337142728327308
$ws.Range("B10").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[8]"
$ws.Range("C10").Value = 36.19047546386719
$ws.Range("D10").Value = 10.85714340209961
$ws.Range("E10").Value = 2.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = "myclk"
$ws.Range("H10").Value = "FF LUT "

# Row 11
$ws.Range("A11").Value = 0.000006589682925550733
$ws.Range("B11").Value = "firConvolutionOperationChaining_IP/U0/shiftRegister_U/firConvolutionOpebkb_ram_U/E[0]"
$ws.Range("C11").Value = 36.507938385009766
$ws.Range("D11").Value = 10.952381134033203
$ws.Range("E11").Value = 5.0
$ws.Range("F11").Value = 3.0
$ws.Range("G11").Value = "myclk"
$ws.Range("H11").Value = "FF LUT "

# Row 12
$ws.Range("A12").Value = 0.000004921904746879591
$ws.Range("B12").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[10]"
$ws.Range("C12").Value = 36.19047546386719
$ws.Range("D12").Value = 10.85714340209961
$ws.Range("E12").Value = 2.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = "myclk"
$ws.Range("H12").Value = "FF LUT "

# Row 13
$ws.Range("A13").Value = 0.000004849364813708235
$ws.Range("B13").Value = "firConvolutionOperationChaining_IP/U0/ap_NS_fsm[5]"
$ws.Range("C13").Value = 43.4920654296875
$ws.Range("D13").Value = 13.047618865966797
$ws.Range("E13").Value = 1.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = "myclk"
$ws.Range("H13").Value = "FF "

# Row 14
$ws.Range("A14").Value = 0.000004704761977336602
$ws.Range("B14").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[7]"
$ws.Range("C14").Value = 36.19047546386719
$ws.Range("D14").Value = 10.85714340209961
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = "myclk"
$ws.Range("H14").Value = "FF LUT "

# Row 15
$ws.Range("A15").Value = 0.000004704761977336602
$ws.Range("B15").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[9]"
$ws.Range("C15").Value = 36.19047546386719
$ws.Range("D15").Value = 10.85714340209961
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = "myclk"
$ws.Range("H15").Value = "FF LUT "

# Row 16
$ws.Range("A16").Value = 0.0000026628572413756046
$ws.Range("B16").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_state5"
$ws.Range("C16").Value = 3.8095240592956543
$ws.Range("D16").Value = 1.1428569555282593
$ws.Range("E16").Value = 35.0
$ws.Range("F16").Value = 9.0
$ws.Range("G16").Value = "myclk"
$ws.Range("H16").Value = "FF LUT "

# Row 17
$ws.Range("A17").Value = 0.0000007276190672200755
$ws.Range("B17").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[2]"
$ws.Range("C17").Value = 3.8095240592956543
$ws.Range("D17").Value = 1.1428569555282593
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 2.0
$ws.Range("G17").Value = "myclk"
$ws.Range("H17").Value = "FF LUT "

# Row 18
$ws.Range("A18").Value = 0.00000037904763416918286
$ws.Range("B18").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[3]"
$ws.Range("C18").Value = 3.8095240592956543
$ws.Range("D18").Value = 1.1428569555282593
$ws.Range("E18").Value = 2.0
$ws.Range("F18").Value = 2.0
$ws.Range("G18").Value = "myclk"
$ws.Range("H18").Value = "FF LUT "

# Row 19
$ws.Range("A19").Value = 0.00000036761903743354196
$ws.Range("B19").Value = "firConvolutionOperationChaining_IP/U0/ap_CS_fsm_reg_n_0_[1]"
$ws.Range("C19").Value = 3.8095240592956543
$ws.Range("D19").Value = 1.1428569555282593
$ws.Range("E19").Value = 2.0
$ws.Range("F19").Value = 2.0
$ws.Range("G19").Value = "myclk"
$ws.Range("H19").Value = "FF LUT "

# Rows 20-22 no longer exist in the updated report; remove them
$ws.Range("A20:H22").Delete()
